$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

$ws.Range("A1").Value = "Columna"
$ws.Range("B1").Value = "Tipo Dato"
$ws.Range("C1").Value = "Tamaño"
$ws.Range("D1").Value = "Descripción"
$ws.Range("A2").Value = "a1"
$ws.Range("B2").Value = "INT"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "ID"
$ws.Range("A3").Value = "a2"
$ws.Range("B3").Value = "INT"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "correlativo"
$ws.Range("E3").Value = "COMPROBANTE DE PAGO"
$ws.Range("A4").Value = "a3"
$ws.Range("B4").Value = "VARCHAR"
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = "Fecha de Emision"
$ws.Range("A5").Value = "a4"
$ws.Range("B5").Value = "VARCHAR"
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "Fecha de Vencimiento"
$ws.Range("A6").Value = "a5"
$ws.Range("B6").Value = "VARCHAR"
$ws.Range("D6").Value = "TIPO (tabla 10)"
$ws.Range("A7").Value = "a6"
$ws.Range("B7").Value = "VARCHAR"
$ws.Range("D7").Value = "Serie"
$ws.Range("A8").Value = "a7"
$ws.Range("B8").Value = "VARCHAR"
$ws.Range("D8").Value = "Numero"
$ws.Range("A9").Value = "a8"
$ws.Range("B9").Value = "VARCHAR"
$ws.Range("D9").Value = "TIPO (tabla 2)"
$ws.Range("E9").Value = "DOCUMENTO DE IDENTIDAD"
$ws.Range("A10").Value = "a9"
$ws.Range("B10").Value = "VARCHAR"
$ws.Range("C10").Value = 11
$ws.Range("D10").Value = "Numero"
$ws.Range("A11").Value = "a10"
$ws.Range("B11").Value = "VARCHAR"
$ws.Range("D11").Value = "Apellidos y nombres"
$ws.Range("A12").Value = "a11"
$ws.Range("B12").Value = "VARCHAR"
$ws.Range("D12").Value = "Cargo"
$ws.Range("E12").Value = "COD. CONTABLE"
$ws.Range("A13").Value = "a12"
$ws.Range("B13").Value = "VARCHAR"
$ws.Range("D13").Value = "Abono"
$ws.Range("A14").Value = "a13"
$ws.Range("B14").Value = "INT"
$ws.Range("C14").Value = 15
$ws.Range("D14").Value = "Valor facturado de la exportacion"
$ws.Range("A15").Value = "a14"
$ws.Range("B15").Value = "INT"
$ws.Range("C15").Value = 15
$ws.Range("D15").Value = "Base imponible de la operación gravada"
$ws.Range("A16").Value = "a15"
$ws.Range("B16").Value = "INT"
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = "Exonerada"
$ws.Range("E16").Value = "Importe de la operación exonerada o inafecta"
$ws.Range("A17").Value = "a16"
$ws.Range("B17").Value = "INT"
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = "Inafecta"
$ws.Range("A18").Value = "a17"
$ws.Range("B18").Value = "INT"
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = "ISC"
$ws.Range("A19").Value = "a18"
$ws.Range("B19").Value = "INT"
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = "IGV o IPM"
$ws.Range("A20").Value = "a19"
$ws.Range("B20").Value = "INT"
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = "Otros tributos y cargos que no forman parte de la base imponible"
$ws.Range("A21").Value = "a20"
$ws.Range("B21").Value = "INT"
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = "Importe total del comprobante de pago"
$ws.Range("A22").Value = "a21"
$ws.Range("B22").Value = "INT"
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = "Tipo de cambio"
$ws.Range("A23").Value = "a22"
$ws.Range("B23").Value = "VARCHAR"
$ws.Range("D23").Value = "COD CTA OPERACIÓN AL CONTADO"
$ws.Range("A24").Value = "a23"
$ws.Range("B24").Value = "VARCHAR"
$ws.Range("D24").Value = "Asiento"

# Header row style (bold white on red, centered, bordered) - same as style index 1
$ws.Range("A1:D1").Style = $wb.Worksheets.Item("db_compras").Range("A1:D1").Style

# Apply center alignment (new style index 4) to column E cells
$ws.Range("E3:E17").HorizontalAlignment = -4108

# Merge cells in column E
$ws.Range("E3:E8").Merge()
$ws.Range("E9:E11").Merge()
$ws.Range("E12:E13").Merge()
$ws.Range("E16:E17").Merge()

# Column widths (best-fit based on content)
$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(5).AutoFit()

# View settings for Hoja2 (sheet4) and db_compras (sheet3)
$ws.Cells.Item(22,3).Select()
$ws.Application.ActiveWindow.ScrollRow = 4

$wsCompras = $wb.Worksheets.Item("db_compras")
$wsCompras.Activate()
$wsCompras.Range("A1:D1").Select()
$wsCompras.Application.ActiveWindow.ScrollRow = 4

# Workbook view: switch to db_compras tab visible and Hoja3 active
$wsHoja3 = $wb.Worksheets.Item("Hoja3")
$wsHoja3.Activate()
